$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "MplusÂ " typo (stray non-breaking-space/mojibake characters) -> "Mplus"
$ws.Range("A36").Value = "Mplus"
$ws.Range("B36").Value = "Mplus"

# Update default font of the workbook's Normal style from Arial to Calibri
$normal = $wb.Styles.Item("Normal")
$normal.Font.Name = "Calibri"

# Reflect the new selection left after editing
$ws.Range("B36").Select() | Out-Null
